$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 14-17 (dataset shrinks from 16 data rows to 12 data rows)
$ws.Range("A14:D17").EntireRow.Delete()

# New data values for rows 2-13 (A:Date serial, B:Channel, C:Metric, D:Value)
# $null entries mean the cell must be fully cleared (no value / no leftover style)
$data = @(
    @(45922, $null,   "Spend", $null),
    @(45922, "TV",    "GRPs",  4),
    @(45908, $null,   "Spend", 88),
    @($null, "TV",    "Spend", 71),
    @(45915, $null,   "GRPs",  $null),
    @(45908, "Radio", "Spend", $null),
    @($null, "TV",    "GRPs",  6),
    @(45922, "Radio", "GRPs",  9),
    @(45915, $null,   $null,   126),
    @(45915, "Radio", "GRPs",  9),
    @(45908, $null,   "GRPs",  2),
    @(45915, "Radio", "Spend", 177)
)

$colLetters = @("A", "B", "C", "D")

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]

    for ($j = 0; $j -lt 4; $j++) {
        $cell = $ws.Range($colLetters[$j] + $row)
        if ($null -eq $vals[$j]) {
            $cell.Clear()
        } else {
            $cell.Value = $vals[$j]
        }
    }
}
